# more bulk upload fixes
# Insert a new "Date Created (Year)*" column after column B (so it becomes
# the new column C), shifting all existing columns from C onward one to
# the right, and fill in a default year value (2000) for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at C; everything from the old C onward shifts right.
$ws.Columns("C:C").Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "Date Created (Year)*"

# Default "Date Created" year values for the three data rows.
$ws.Range("C2").Value = 2000
$ws.Range("C3").Value = 2000
$ws.Range("C4").Value = 2000

# Give the new values an explicit black font color (new style).
$ws.Range("C2:C4").Font.Color = 0

# Update the active selection to the newly inserted column.
$ws.Range("C1:C4").Select()
